# Auto-generated edit script: update cached Leve profit-calc values
# per the scheduled-runner refresh (Sheets diff).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 891.44446
$ws.Cells.Item(18, 9).Value = 1393.1666
$ws.Cells.Item(18, 10).Value = 640.5833
$ws.Cells.Item(18, 11).Value = 1393.1666
$ws.Cells.Item(18, 12).Value = 640.5833
$ws.Cells.Item(18, 13).Value = -1109.1666
$ws.Cells.Item(18, 14).Value = -1208.5833

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 4216.3335
$ws.Cells.Item(64, 9).Value = 3634.423
$ws.Cells.Item(64, 11).Value = 3634.423
$ws.Cells.Item(64, 13).Value = -3386.423

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(67, 8).Value = 4216.3335
$ws.Cells.Item(67, 9).Value = 3634.423
$ws.Cells.Item(67, 11).Value = 3634.423
$ws.Cells.Item(67, 13).Value = -2776.423

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(106, 8).Value = 7272.5
$ws.Cells.Item(106, 9).Value = 7272.5
$ws.Cells.Item(106, 11).Value = 7272.5
$ws.Cells.Item(106, 13).Value = -6641.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 16620.117
$ws.Cells.Item(116, 9).Value = 4308.636
$ws.Cells.Item(116, 10).Value = 39191.168
$ws.Cells.Item(116, 11).Value = 4308.636
$ws.Cells.Item(116, 12).Value = 39191.168
$ws.Cells.Item(116, 13).Value = -866.6360000000004
$ws.Cells.Item(116, 14).Value = -46075.168

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 5540.904
$ws.Cells.Item(132, 9).Value = 5363.0835
$ws.Cells.Item(132, 11).Value = 16089.2505
$ws.Cells.Item(132, 13).Value = -13559.2505

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 3164.0667
$ws.Cells.Item(138, 9).Value = 2468.818
$ws.Cells.Item(138, 10).Value = 3566.5789
$ws.Cells.Item(138, 11).Value = 7406.454000000001
$ws.Cells.Item(138, 12).Value = 10699.7367
$ws.Cells.Item(138, 13).Value = -2266.454000000001
$ws.Cells.Item(138, 14).Value = -20979.7367

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 144427.16
$ws.Cells.Item(32, 9).Value = 213584.12
$ws.Cells.Item(32, 11).Value = 213584.12
$ws.Cells.Item(32, 13).Value = -213297.12

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1836.8667
$ws.Cells.Item(45, 9).Value = 1755.3
$ws.Cells.Item(45, 11).Value = 1755.3
$ws.Cells.Item(45, 13).Value = -1378.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2504399
$ws.Cells.Item(61, 9).Value = 4420.4116
$ws.Cells.Item(61, 11).Value = 4420.4116
$ws.Cells.Item(61, 13).Value = -4208.4116

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 628187.0600000001
$ws.Cells.Item(132, 9).Value = 644189.3
$ws.Cells.Item(132, 11).Value = 1932567.9
$ws.Cells.Item(132, 13).Value = -1930037.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(134, 8).Value = 50248.625
$ws.Cells.Item(134, 10).Value = 50248.625
$ws.Cells.Item(134, 12).Value = 50248.625
$ws.Cells.Item(134, 14).Value = -60388.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 2504399
$ws.Cells.Item(136, 9).Value = 4420.4116
$ws.Cells.Item(136, 11).Value = 13261.2348
$ws.Cells.Item(136, 13).Value = -10711.2348

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(11, 8).Value = 68.8
$ws.Cells.Item(11, 10).Value = 49
$ws.Cells.Item(11, 12).Value = 49
$ws.Cells.Item(11, 14).Value = -329

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(96, 8).Value = 22147.46
$ws.Cells.Item(96, 9).Value = 18993.084
$ws.Cells.Item(96, 11).Value = 18993.084
$ws.Cells.Item(96, 13).Value = -16247.084

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 8389.852000000001
$ws.Cells.Item(105, 9).Value = 7791.4
$ws.Cells.Item(105, 10).Value = 10099.714
$ws.Cells.Item(105, 11).Value = 7791.4
$ws.Cells.Item(105, 12).Value = 10099.714
$ws.Cells.Item(105, 13).Value = -6044.4
$ws.Cells.Item(105, 14).Value = -13593.714

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1276.6666
$ws.Cells.Item(107, 9).Value = 975
$ws.Cells.Item(107, 10).Value = 1880
$ws.Cells.Item(107, 11).Value = 975
$ws.Cells.Item(107, 12).Value = 1880
$ws.Cells.Item(107, 13).Value = 945
$ws.Cells.Item(107, 14).Value = -5720

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2457204.5
$ws.Cells.Item(134, 9).Value = 4558.54
$ws.Cells.Item(134, 10).Value = 9270110
$ws.Cells.Item(134, 11).Value = 13675.62
$ws.Cells.Item(134, 12).Value = 27810330
$ws.Cells.Item(134, 13).Value = -11140.62
$ws.Cells.Item(134, 14).Value = -27815400

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1795106.6
$ws.Cells.Item(31, 9).Value = 1986871.6
$ws.Cells.Item(31, 10).Value = 5299.6665
$ws.Cells.Item(31, 11).Value = 1986871.6
$ws.Cells.Item(31, 12).Value = 5299.6665
$ws.Cells.Item(31, 13).Value = -1986576.6
$ws.Cells.Item(31, 14).Value = -5889.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 1795106.6
$ws.Cells.Item(34, 9).Value = 1986871.6
$ws.Cells.Item(34, 10).Value = 5299.6665
$ws.Cells.Item(34, 11).Value = 1986871.6
$ws.Cells.Item(34, 12).Value = 5299.6665
$ws.Cells.Item(34, 13).Value = -1986669.6
$ws.Cells.Item(34, 14).Value = -5703.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 14924.5
$ws.Cells.Item(86, 9).Value = 18267.715
$ws.Cells.Item(86, 10).Value = 7123.6665
$ws.Cells.Item(86, 11).Value = 18267.715
$ws.Cells.Item(86, 12).Value = 7123.6665
$ws.Cells.Item(86, 13).Value = -17144.715
$ws.Cells.Item(86, 14).Value = -9369.666499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(89, 8).Value = 14924.5
$ws.Cells.Item(89, 9).Value = 18267.715
$ws.Cells.Item(89, 10).Value = 7123.6665
$ws.Cells.Item(89, 11).Value = 91338.575
$ws.Cells.Item(89, 12).Value = 35618.3325
$ws.Cells.Item(89, 13).Value = -85722.575
$ws.Cells.Item(89, 14).Value = -46850.3325

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 6994.2856
$ws.Cells.Item(105, 9).Value = 8136.143
$ws.Cells.Item(105, 10).Value = 4710.5713
$ws.Cells.Item(105, 11).Value = 8136.143
$ws.Cells.Item(105, 12).Value = 4710.5713
$ws.Cells.Item(105, 13).Value = -6389.143
$ws.Cells.Item(105, 14).Value = -8204.5713

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 7522.2144
$ws.Cells.Item(122, 9).Value = 1689.7037
$ws.Cells.Item(122, 11).Value = 5069.1111
$ws.Cells.Item(122, 13).Value = -2619.1111

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 3277820.2
$ws.Cells.Item(5, 9).Value = 5953365.5
$ws.Cells.Item(5, 10).Value = 1940047.6
$ws.Cells.Item(5, 11).Value = 17860096.5
$ws.Cells.Item(5, 12).Value = 5820142.800000001
$ws.Cells.Item(5, 13).Value = -17859984.5
$ws.Cells.Item(5, 14).Value = -5820366.800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 2651.6667
$ws.Cells.Item(113, 10).Value = 2814.1538
$ws.Cells.Item(113, 12).Value = 8442.4614
$ws.Cells.Item(113, 14).Value = -12782.4614

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(135, 8).Value = 3277820.2
$ws.Cells.Item(135, 9).Value = 5953365.5
$ws.Cells.Item(135, 10).Value = 1940047.6
$ws.Cells.Item(135, 11).Value = 53580289.5
$ws.Cells.Item(135, 12).Value = 17460428.4
$ws.Cells.Item(135, 13).Value = -53577754.5
$ws.Cells.Item(135, 14).Value = -17465498.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(136, 8).Value = 8904.166999999999
$ws.Cells.Item(136, 10).Value = 14874.75
$ws.Cells.Item(136, 12).Value = 44624.25
$ws.Cells.Item(136, 14).Value = -54824.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3237.8572
$ws.Cells.Item(80, 9).Value = 2708
$ws.Cells.Item(80, 11).Value = 2708
$ws.Cells.Item(80, 13).Value = -1710

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(82, 8).Value = 0
$ws.Cells.Item(82, 9).Value = 0
$ws.Cells.Item(82, 11).Value = 0
$ws.Cells.Item(82, 13).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 3237.8572
$ws.Cells.Item(83, 9).Value = 2708
$ws.Cells.Item(83, 11).Value = 13540
$ws.Cells.Item(83, 13).Value = -8548

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(85, 8).Value = 0
$ws.Cells.Item(85, 9).Value = 0
$ws.Cells.Item(85, 11).Value = 0
$ws.Cells.Item(85, 13).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 2748.7058
$ws.Cells.Item(113, 9).Value = 2111.818
$ws.Cells.Item(113, 11).Value = 2111.818
$ws.Cells.Item(113, 13).Value = 58.18199999999979

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 5241.2104
$ws.Cells.Item(122, 10).Value = 5399.75
$ws.Cells.Item(122, 12).Value = 16199.25
$ws.Cells.Item(122, 14).Value = -21099.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 10262.681
$ws.Cells.Item(132, 9).Value = 8499.902
$ws.Cells.Item(132, 11).Value = 25499.706
$ws.Cells.Item(132, 13).Value = -22969.706

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 5240.8237
$ws.Cells.Item(46, 9).Value = 1999
$ws.Cells.Item(46, 10).Value = 5443.4375
$ws.Cells.Item(46, 11).Value = 1999
$ws.Cells.Item(46, 12).Value = 5443.4375
$ws.Cells.Item(46, 13).Value = -1811
$ws.Cells.Item(46, 14).Value = -5819.4375

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 10669.526
$ws.Cells.Item(61, 9).Value = 12259.3125
$ws.Cells.Item(61, 11).Value = 12259.3125
$ws.Cells.Item(61, 13).Value = -12057.3125

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 10669.526
$ws.Cells.Item(113, 9).Value = 12259.3125
$ws.Cells.Item(113, 11).Value = 12259.3125
$ws.Cells.Item(113, 13).Value = -10089.3125

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 3287.5881
$ws.Cells.Item(122, 9).Value = 3081.3635
$ws.Cells.Item(122, 11).Value = 9244.0905
$ws.Cells.Item(122, 13).Value = -6794.0905

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 3654435.5
$ws.Cells.Item(132, 10).Value = 3584.5454
$ws.Cells.Item(132, 12).Value = 10753.6362
$ws.Cells.Item(132, 14).Value = -15813.6362

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 10424774
$ws.Cells.Item(136, 9).Value = 12503436
$ws.Cells.Item(136, 11).Value = 37510308
$ws.Cells.Item(136, 13).Value = -37507758

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 9).Value = 1449.4286
$ws.Cells.Item(96, 10).Value = 1633.3334
$ws.Cells.Item(96, 11).Value = 1449.4286
$ws.Cells.Item(96, 12).Value = 1633.3334
$ws.Cells.Item(96, 13).Value = -76.42859999999996
$ws.Cells.Item(96, 14).Value = -4379.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 43116.43
$ws.Cells.Item(122, 9).Value = 3471.3684
$ws.Cells.Item(122, 11).Value = 10414.1052
$ws.Cells.Item(122, 13).Value = -7964.1052

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 6412209.5
$ws.Cells.Item(132, 9).Value = 7577511.5
$ws.Cells.Item(132, 10).Value = 3049.75
$ws.Cells.Item(132, 11).Value = 22732534.5
$ws.Cells.Item(132, 12).Value = 9149.25
$ws.Cells.Item(132, 13).Value = -22730004.5
$ws.Cells.Item(132, 14).Value = -14209.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 13527520
$ws.Cells.Item(136, 9).Value = 2718397.2
$ws.Cells.Item(136, 11).Value = 8155191.600000001
$ws.Cells.Item(136, 13).Value = -8152641.600000001
